{"js": "// Apply the benchmark table updates described by the commit:\n// \"Fixed README.md stats and docx preparation for all DaCapo - JDK 17 -\n//  Shenandoah GC tests\"\n//\n// The document is a single-column table; each row holds one stat value.\n// We update specific rows (by 0-based row index) to their corrected\n// values, and collapse the three \"raw metrics dump\" rows (which held a\n// tab-separated run of numbers) down to the single summary value they\n// should have shown.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// rowIndex -> new cell text\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"124\",\n  5: \"0.00070\",\n  6: \"0.00019\",\n  7: \"0.00008\",\n  8: \"0.00025\",\n  9: \"0.00032\",\n  10: \"0.00042\",\n  11: \"0.02407\",\n  43: \"100\",\n  44: \"0.02\",\n  45: \"2364\",\n};\n\nfor (const [rowIndex, newText] of Object.entries(updates)) {\n  const cell = table.getCell(Number(rowIndex), 0);\n  cell.value = newText;\n}\n\nawait context.sync();\n", "ps1": "# Apply the benchmark table updates described by the commit:\n# \"Fixed README.md stats and docx preparation for all DaCapo - JDK 17 -\n#  Shenandoah GC tests\"\n#\n# The document is a single-column table; each row holds one stat value.\n# We update specific rows (by 1-based Word COM row index) to their\n# corrected values, and collapse the three \"raw metrics dump\" rows\n# (which held a tab-separated run of numbers) down to the single\n# summary value they should have shown.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$updates = @{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"124\"\n    6  = \"0.00070\"\n    7  = \"0.00019\"\n    8  = \"0.00008\"\n    9  = \"0.00025\"\n    10 = \"0.00032\"\n    11 = \"0.00042\"\n    12 = \"0.02407\"\n    44 = \"100\"\n    45 = \"0.02\"\n    46 = \"2364\"\n}\n\nforeach ($rowIndex in $updates.Keys) {\n    $cell = $t.Cell($rowIndex, 1)\n    $cell.Range.Text = $updates[$rowIndex]\n}\n"}
